# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Piña"
# This shifts existing rows 150-156 down to 151-157, and populates the new
# row 150 with the new data point (week of 2021-11-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 150; this pushes rows 150-156 down to 151-157
# and carries their formatting/content along.
$ws.Rows.Item(150).Insert()

# Copy the style (formatting) of the date cell from the row above (D149)
# into the new D150 cell so it keeps the date number format.
$ws.Range("D149").Copy()
$ws.Range("D150").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Populate the new row 150 with values matching the rest of the table.
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44509
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100108
$ws.Range("H150").Value = "Tropicales y subtropicales"
$ws.Range("I150").Value = 100108005
$ws.Range("J150").Value = "Piña"
$ws.Range("K150").Value = "Caramelo"
$ws.Range("L150").Value = "Segunda"
$ws.Range("M150").Value = 200
$ws.Range("N150").Value = 22000
$ws.Range("O150").Value = 23000
$ws.Range("P150").Value = 22500
$ws.Range("Q150").Value = "$/caja 14 unidades"
$ws.Range("R150").Value = "Ecuador"
$ws.Range("S150").Value = 1607
$ws.Range("T150").Value = 14

$wb.Save()
